$d = $word.ActiveDocument

function Replace-WholeWord($find, $replace) {
    # MatchWholeWord=$true keeps this from matching inside longer tokens
    # (e.g. "L" inside "XL", or "M" inside "HERMAWAN"/"FIRMANDA").
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# NO fields (each value is unique and alone in its own run/paragraph)
Replace-WholeWord "F39" "F42"
Replace-WholeWord "G1" "F43"
Replace-WholeWord "G2" "F44"

# NAMA fields
Replace-Text "IWAN HERMAWAN" "PARJO DWI KURNIAWAN"
Replace-Text "FIRMANDA DWI KURNIAWAN" "DONA ANDHIKA FAIZAL"
Replace-Text "DWI PRASETYO" "MILYATER H HALOHO"

# KELAS fields: both occurrences of "DP 1 TEKNIKA / 9" change identically to "DP 1 NAUTIKA / 9"
Replace-Text "DP 1 TEKNIKA / 9" "DP 1 NAUTIKA / 9"

# KAOS sizes: column2 L -> XL and column3 XL -> L is a swap, so stage column3's XL
# through a placeholder first so it can't collide with column2's new XL value.
# Column1's M -> XL is done last so it can't collide with the swap either.
Replace-WholeWord "XL" "ZZTMPZZ"
Replace-WholeWord "L" "XL"
Replace-WholeWord "ZZTMPZZ" "L"
Replace-WholeWord "M" "XL"
